# 22 April 1st update
# Inserts two new date columns (21/04/2020 and 22/04/2020) into the wide
# COVID state table, between 21/03/2020 (AQ) and 23/03/2020 (old AR).
# Also corrects Jharkhand's 20/04/2020 figure (1 -> 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new date columns. Inserting at AR first pushes the old
# AR..BB block one column right (old AR "22/03/2020" -> AS). Inserting at
# AT next pushes that shifted block right again, leaving AR and AT as two
# fresh, blank columns for the new dates.
$ws.Columns("AR").Insert()
$ws.Columns("AT").Insert()

# New column headers
$ws.Range("AR1").Value = "21/04/2020"
$ws.Range("AT1").Value = "22/04/2020"

# New data for 21/04/2020 (column AR), by state row
$ws.Range("AR3").Value = 1    # Andaman and Nicobar Islands
$ws.Range("AR4").Value = 35   # Andhra Pradesh
$ws.Range("AR7").Value = 13   # Bihar
$ws.Range("AR8").Value = 1    # Chandigarh
$ws.Range("AR10").Value = 75  # Delhi
$ws.Range("AR12").Value = 239 # Gujarat
$ws.Range("AR13").Value = 4   # Haryana
$ws.Range("AR15").Value = 12  # Jammu and Kashmir
$ws.Range("AR17").Value = 10  # Karnataka
$ws.Range("AR18").Value = 19  # Kerala
$ws.Range("AR20").Value = 67  # Madhya Pradesh
$ws.Range("AR21").Value = 552 # Maharashtra
$ws.Range("AR23").Value = 1   # Meghalaya
$ws.Range("AR25").Value = 5   # Odisha
$ws.Range("AR27").Value = 6   # Punjab
$ws.Range("AR28").Value = 159 # Rajasthan
$ws.Range("AR29").Value = 76  # Tamil Nadu
$ws.Range("AR30").Value = 56  # Telangana
$ws.Range("AR32").Value = 153 # Uttar Pradesh
$ws.Range("AR34").Value = 53  # West Bengal

# New data for 22/04/2020 (column AT), by state row
$ws.Range("AT25").Value = 3   # Odisha
$ws.Range("AT34").Value = 31  # West Bengal

# Data correction: Jharkhand's 20/04/2020 figure (column AP, row 16)
$ws.Range("AP16").Value = 5
